# The experimental data row for bin B=6,C=7 (pT 4.4-4.6) had a negative
# cross-section value (I44 = -1.5E-3, K44 = -3.75E-4), which is not
# physically meaningful. Remove that entire row, shifting all following
# rows up by one (matching the commit message "remove negative cross
# sections").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(44).Delete()
